$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Column D prices that are plain numeric-looking strings are written with a
# leading apostrophe so Excel keeps them as text (matching the source data,
# which stores every price as text) instead of auto-converting to a number.

$ws.Range("D2").Value = "66.189.16"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.168.43"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'606.22"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'154.02"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.168.42"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "'0.546"
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -9.09%  "
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("D13").Value = "'0.0000266"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").Value = "'38.29"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "3.688.49"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").Value = "66.189.99"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'7.40"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "3.172.68"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "'510.32"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'15.38"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "'7.99"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "'14.77"
$ws.Range("E24").Value = "  -3.58%  "
$ws.Range("D25").Value = "'84.62"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "'9.09"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "'2.37"
$ws.Range("E29").Value = "  +5.25%  "
$ws.Range("D30").Value = "'3.03"
$ws.Range("E30").Value = "  +6.29%  "
$ws.Range("D31").Value = "'7.17"
$ws.Range("E31").Value = "  +4.77%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("D35").Value = "'6.51"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").Value = "'499.86"
$ws.Range("E36").Value = "  +4.73%  "
$ws.Range("D37").Value = "'54.90"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").Value = "'0.0880"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  +7.37%  "
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("D42").Value = "0.0₃0683"
$ws.Range("E42").Value = "  +6.47%  "

# Rows 43/44 swapped ranking order: dogwifhat now above TheGraph.
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.81"
$ws.Range("E43").Value = "  -4.60%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.296"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'2.42"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").Value = "2.824.54"
$ws.Range("E46").Value = "  -4.34%  "

# Rows 48/49 swapped ranking order: USDe now above ThetaToken.
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "'2.37"
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "'35.26"
$ws.Range("E51").Value = "  +6.09%  "
